$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 137, shifting existing rows 137:184 down to 138:185.
$ws.Rows("137:137").Insert()

# Populate the newly inserted row 137 with the new weekly data point.
$ws.Cells.Item(137, 1).Value = 6
$ws.Cells.Item(137, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(137, 3).Value = "Metropolitana"
$ws.Cells.Item(137, 4).Value = 44636
$ws.Cells.Item(137, 5).Value = 13
$ws.Cells.Item(137, 6).Value = "Fruta"
$ws.Cells.Item(137, 7).Value = 100101
$ws.Cells.Item(137, 8).Value = "Berries"
$ws.Cells.Item(137, 9).Value = 100101004
$ws.Cells.Item(137, 10).Value = "Frambuesa"
$ws.Cells.Item(137, 11).Value = "Sin especificar"
$ws.Cells.Item(137, 12).Value = "Especial"
$ws.Cells.Item(137, 13).Value = 500
$ws.Cells.Item(137, 14).Value = 8000
$ws.Cells.Item(137, 15).Value = 8000
$ws.Cells.Item(137, 16).Value = 8000
$ws.Cells.Item(137, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(137, 18).Value = "Provincia de Linares"
$ws.Cells.Item(137, 19).Value = 4000
$ws.Cells.Item(137, 20).Value = 2
